$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$ws.Range("D2").Value = "60.448.80"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "2.599.39"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  -0.21%  "

$helper.Formula = "=""513.88"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.69%  "

$helper.Formula = "=""153.34"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("E7").Value = "  +0.20%  "

$helper.Formula = "=""0.598"""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +3.32%  "

$helper.Formula = "=""6.64"""
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("D13").Value = "3.054.85"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").Value = "60.467.75"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "2.618.18"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("E18").Value = "  -1.24%  "

$helper.Formula = "=""358.53"""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +3.84%  "

$helper.Formula = "=""10.56"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +1.67%  "

$helper.Formula = "=""6.19"""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.78%  "

$ws.Range("E22").Value = "  +0.19%  "

$helper.Formula = "=""61.04"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.88%  "

$helper.Formula = "=""0.426"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.91%  "

$ws.Range("D25").Value = "2.716.09"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("E26").Value = "  +0.20%  "

$helper.Formula = "=""0.994"""
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").Value = "0.0₃0834"
$ws.Range("E28").Value = "  -2.55%  "

$ws.Range("E29").Value = "  -2.95%  "

$ws.Range("E30").Value = "  +0.15%  "

$helper.Formula = "=""19.40"""
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("E32").Value = "  +1.64%  "

$helper.Formula = "=""5.94"""
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +3.31%  "

$helper.Formula = "=""150.72"""
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.42%  "

$helper.Formula = "=""4.01"""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$helper.Formula = "=""0.915"""
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +6.66%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$helper.Formula = "=""1.19"""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  +0.12%  "

$helper.Formula = "=""36.30"""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +1.66%  "

$helper.Formula = "=""0.842"""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.83%  "

$helper.Formula = "=""3.75"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.32%  "

$helper.Formula = "=""287.39"""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E45").Value = "  +0.17%  "

$helper.Formula = "=""0.0554"""
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -2.62%  "

$helper.Formula = "=""19.57"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.65%  "

$helper.Formula = "=""4.95"""
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -0.56%  "

$helper.Formula = "=""0.0235"""
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("E50").Value = "  +0.42%  "

$helper.Formula = "=""19.18"""
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +9.26%  "

$helper.Clear()
$excel.CutCopyMode = $false
